$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1032827490704891
$ws.Range("B3").Value = 0.004763625017518802
$ws.Range("C3").Value = 0.0008438719170554776
$ws.Range("D3").Value = 5.481658796936038
$ws.Range("E3").Value = 0.1337365934067819
$ws.Range("F3").Value = 0.003109660717741514
$ws.Range("G3").Value = 0.006417589317296092
$ws.Range("H3").Value = 0.1080463740880079
$ws.Range("B4").Value = 0.009274788942493036
$ws.Range("C4").Value = 0.00133149388849625
$ws.Range("D4").Value = 9.993172196100661
$ws.Range("E4").Value = 0.04521872179843062
$ws.Range("F4").Value = 0.006665100332934794
$ws.Range("G4").Value = 0.01188447755205128
$ws.Range("H4").Value = 0.1125575380129821
$ws.Range("B5").Value = 0.007588028959536758
$ws.Range("C5").Value = 0.005155505738442479
$ws.Range("D5").Value = 4.816190498694034
$ws.Range("E5").Value = 0.07189748726389381
$ws.Range("F5").Value = -0.00251661501146993
$ws.Range("G5").Value = 0.01769267293054345
$ws.Range("H5").Value = 0.1108707780300259
$ws.Range("B6").Value = 0.001104510828489702
$ws.Range("C6").Value = 0.003024503769602311
$ws.Range("D6").Value = 2.285707903482943
$ws.Range("E6").Value = 0.005270529431042424
$ws.Range("F6").Value = -0.004823425188889096
$ws.Range("G6").Value = 0.007032446845868501
$ws.Range("H6").Value = 0.1043872598989788
$ws.Range("B7").Value = 0.005254766853219207
$ws.Range("C7").Value = 0.004644266476306593
$ws.Range("D7").Value = 4.379847309338646
$ws.Range("E7").Value = 0.001643401456486382
$ws.Range("F7").Value = -0.00384785499573886
$ws.Range("G7").Value = 0.01435738870217728
$ws.Range("H7").Value = 0.1085375159237083
$ws.Range("B8").Value = 0.02934553417950313
$ws.Range("C8").Value = 0.00425728056543912
$ws.Range("D8").Value = 11.36188312237061
$ws.Range("E8").Value = 0.04708332184117783
$ws.Range("F8").Value = 0.0210013939927424
$ws.Range("G8").Value = 0.03768967436626387
$ws.Range("H8").Value = 0.1326282832499922
$ws.Range("B9").Value = 0.03265094753321694
$ws.Range("C9").Value = 0.005372996153359723
$ws.Range("D9").Value = 15.57276322578737
$ws.Range("E9").Value = 0.06189039294776968
$ws.Range("F9").Value = 0.02212003928590978
$ws.Range("G9").Value = 0.04318185578052409
$ws.Range("H9").Value = 0.1359336966037061
$ws.Range("B10").Value = -0.1032827490704891
$ws.Range("C10").Value = 0.0005694284115987781
$ws.Range("D10").Value = -230.4098665142313
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.1043988120543223
$ws.Range("G10").Value = -0.1021666860866559
$ws.Range("B11").Value = -0.04736442151549015
$ws.Range("C11").Value = 0.0006016225185623984
$ws.Range("D11").Value = -96.1312663314131
$ws.Range("E11").Value = [double]"1.318191952768195e-279"
$ws.Range("F11").Value = -0.04854358400135905
$ws.Range("G11").Value = -0.04618525902962124
$ws.Range("H11").Value = 0.05591832755499895
$ws.Range("B12").Value = -0.03909081751746011
$ws.Range("C12").Value = 0.0005869529006753061
$ws.Range("D12").Value = -82.57571264211641
$ws.Range("E12").Value = [double]"3.047563889247479e-227"
$ws.Range("F12").Value = -0.04024122798973506
$ws.Range("G12").Value = -0.03794040704518515
$ws.Range("H12").Value = 0.06419193155302899
$ws.Range("B13").Value = -0.03335510558858117
$ws.Range("C13").Value = 0.0005774692866451155
$ws.Range("D13").Value = -71.948214187369
$ws.Range("E13").Value = [double]"6.632606906660757e-126"
$ws.Range("F13").Value = -0.0344869284627916
$ws.Range("G13").Value = -0.03222328271437074
$ws.Range("H13").Value = 0.06992764348190794
$ws.Range("B14").Value = -0.03087808215686661
$ws.Range("C14").Value = 0.0005624301923174425
$ws.Range("D14").Value = -69.44666938940283
$ws.Range("E14").Value = [double]"4.98664466004938e-140"
$ws.Range("F14").Value = -0.0319804288462914
$ws.Range("G14").Value = -0.02977573546744182
$ws.Range("H14").Value = 0.07240466691362249
$ws.Range("B15").Value = -0.02650656905974687
$ws.Range("C15").Value = 0.0005591936326301471
$ws.Range("D15").Value = -58.93009461994982
$ws.Range("E15").Value = [double]"5.556829096270904e-99"
$ws.Range("F15").Value = -0.02760257218820355
$ws.Range("G15").Value = -0.02541056593129019
$ws.Range("H15").Value = 0.07677618001074224
$ws.Range("B16").Value = -0.02414408901062102
$ws.Range("C16").Value = 0.000559660158224322
$ws.Range("D16").Value = -53.43715554716025
$ws.Range("E16").Value = [double]"1.429972428663068e-55"
$ws.Range("F16").Value = -0.02524100652471994
$ws.Range("G16").Value = -0.02304717149652211
$ws.Range("H16").Value = 0.07913866005986808
$ws.Range("B17").Value = -0.02248066792742855
$ws.Range("C17").Value = 0.0005655428809024975
$ws.Range("D17").Value = -49.67194242402783
$ws.Range("E17").Value = [double]"1.654191259525222e-59"
$ws.Range("F17").Value = -0.02358911539882896
$ws.Range("G17").Value = -0.02137222045602814
$ws.Range("H17").Value = 0.08080208114306056
$ws.Range("B18").Value = -0.01989744723870494
$ws.Range("C18").Value = 0.0005687573121337462
$ws.Range("D18").Value = -40.94868575061944
$ws.Range("E18").Value = [double]"7.88864485388906e-43"
$ws.Range("F18").Value = -0.02101219490568048
$ws.Range("G18").Value = -0.0187826995717294
$ws.Range("H18").Value = 0.08338530183178416
$ws.Range("B19").Value = -0.01685928269530315
$ws.Range("C19").Value = 0.0005721717598939574
$ws.Range("D19").Value = -33.2383387244386
$ws.Range("E19").Value = [double]"7.577811518708963e-17"
$ws.Range("F19").Value = -0.01798072257819874
$ws.Range("G19").Value = -0.01573784281240755
$ws.Range("H19").Value = 0.08642346637518596
$ws.Range("B20").Value = -0.01479964664143699
$ws.Range("C20").Value = 0.0005781644270461732
$ws.Range("D20").Value = -28.65157037572939
$ws.Range("E20").Value = 0.01159313765412665
$ws.Range("F20").Value = -0.01593283198543198
$ws.Range("G20").Value = -0.01366646129744199
$ws.Range("H20").Value = 0.08848310242905212
$ws.Range("B21").Value = -0.01086289039593508
$ws.Range("C21").Value = 0.0005905602168530004
$ws.Range("D21").Value = -19.43462176241935
$ws.Range("E21").Value = [double]"3.430133750637629e-07"
$ws.Range("F21").Value = -0.01202037113353169
$ws.Range("G21").Value = -0.009705409658338469
$ws.Range("H21").Value = 0.09241985867455403
$ws.Range("B22").Value = -0.007222199894794783
$ws.Range("C22").Value = 0.0005930676062191558
$ws.Range("D22").Value = -11.12689255830032
$ws.Range("E22").Value = 0.04741263110994309
$ws.Range("F22").Value = -0.008384595047250537
$ws.Range("G22").Value = -0.006059804742339027
$ws.Range("H22").Value = 0.09606054917569433
$ws.Range("B23").Value = -0.005710960279529289
$ws.Range("C23").Value = 0.0005985555987679838
$ws.Range("D23").Value = -8.569398641180443
$ws.Range("E23").Value = 0.1364817949179106
$ws.Range("F23").Value = -0.006884111729868276
$ws.Range("G23").Value = -0.004537808829190301
$ws.Range("H23").Value = 0.09757178879095982
$ws.Range("B24").Value = -0.004884883832734936
$ws.Range("C24").Value = 0.0005875874666679092
$ws.Range("D24").Value = -8.585938013147771
$ws.Range("E24").Value = 0.05745343707799912
$ws.Range("F24").Value = -0.006036538052667169
$ws.Range("G24").Value = -0.003733229612802705
$ws.Range("H24").Value = 0.09839786523775418
$ws.Range("B25").Value = -0.00284536375971885
$ws.Range("C25").Value = 0.0005939617381392074
$ws.Range("D25").Value = -5.770895618908816
$ws.Range("E25").Value = 0.008529766938610656
$ws.Range("F25").Value = -0.00400951135889401
$ws.Range("G25").Value = -0.00168121616054369
$ws.Range("H25").Value = 0.1004373853107703
$ws.Range("B26").Value = 0.02723554208481277
$ws.Range("C26").Value = 0.0007842127250260329
$ws.Range("D26").Value = 25.86458914183782
$ws.Range("E26").Value = [double]"6.906694836983625e-15"
$ws.Range("F26").Value = 0.02569850850163167
$ws.Range("G26").Value = 0.02877257566799387
$ws.Range("H26").Value = 0.1305182911553019
